$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.122.34"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "2.551.39"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'585.08"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").Value = "'147.45"
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.585"
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("E10").Value = "  -3.92%  "
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "'0.356"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("E13").Value = "  -3.83%  "
$ws.Range("D14").Value = "3.001.74"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "62.979.37"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "2.550.72"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "'11.37"
$ws.Range("E18").Value = "  -2.74%  "
$ws.Range("D19").Value = "'337.36"
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("D20").Value = "'4.34"
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("D21").Value = "'6.79"
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Value = "'65.94"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").Value = "'8.41"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").Value = "'7.71"
$ws.Range("E29").Value = "  +8.49%  "
$ws.Range("E30").Value = "  +5.21%  "
$ws.Range("D31").Value = "0.0₃0818"
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("D32").Value = "'178.24"
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'1.55"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value = "'417.85"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("D35").Value = "'19.20"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "'0.401"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D38").Value = "'4.37"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'39.65"
$ws.Range("E41").Value = "  -1.01%  "
$ws.Range("D42").Value = "'150.78"
$ws.Range("E42").Value = "  -2.91%  "
$ws.Range("D43").Value = "'3.80"
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").Value = "'20.87"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("D45").Value = "'0.0543"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.604"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.0973"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").Value = "'0.0239"
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("D49").Value = "'18.33"
$ws.Range("E49").Value = "  -1.97%  "
$ws.Range("E50").Value = "  -5.76%  "
$ws.Range("E51").Value = "  -0.41%  "
